$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the amounts in column C (rows 2-11)
$ws.Range("C2").Value = 1000
$ws.Range("C3").Value = 3000
$ws.Range("C4").Value = 1000
$ws.Range("C5").Value = 1000
$ws.Range("C6").Value = 1000
$ws.Range("C7").Value = 1000
$ws.Range("C8").Value = 1000
$ws.Range("C9").Value = 1000
$ws.Range("C10").Value = 1000
$ws.Range("C11").Value = 1000

# Move the selection down to C4:C11 (active cell C4), matching the saved view state
$ws.Range("C4:C11").Select()
